# Fix calculation of correct Excel column width / comparison of column
# widths when columns belong to a column group.
#
# Net effect captured by this change for each worksheet:
#   - Row 1 gets an explicit height of 24pt (ht="24" customHeight="1").
#   - The stored column widths for columns A, E and F are a tiny
#     (sub-pixel) floating point correction of the already-stored width
#     (e.g. 5.5709375 -> 5.5703125, 24.5709375 -> 24.5703125). These are
#     not perceptible width changes (same integer pixel width), just a
#     more precise internal float, so we leave the already-correct
#     widths as-is rather than risk moving them further away by round
#     tripping them through the (coarser, pixel-snapping) ColumnWidth
#     setter.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Rows.Item(1).RowHeight = 24
}
